$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.593.84"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.594.24"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'210.86"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'0.514"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "'19.42"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.819.50"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.597.01"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "'64.40"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "26.587.32"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'207.70"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'145.21"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'15.22"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.651"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "1.279.92"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'63.69"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "'0.918"
$ws.Range("E45").Value = "  +9.57%  "
$ws.Range("D46").Value = "1.731.02"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'89.49"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("E50").Value = "  +4.10%  "
$ws.Range("E51").Value = "  +0.82%  "
